$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 449.0375  # H17: was 459.86517
$ws.Cells.Item(17, 10).Value = 449.0375  # J17: was 459.86517
$ws.Cells.Item(17, 12).Value = 1347.1125  # L17: was 1379.59551
$ws.Cells.Item(17, 14).Value = -1683.1125  # N17: was -1715.59551
$ws.Cells.Item(33, 8).Value = 442.05884  # H33: was 461
$ws.Cells.Item(33, 10).Value = 569.5  # J33: was 1000
$ws.Cells.Item(33, 12).Value = 569.5  # L33: was 1000
$ws.Cells.Item(33, 14).Value = -1027.5  # N33: was -1458
$ws.Cells.Item(40, 8).Value = 1835.6428  # H40: was 1764.2142
$ws.Cells.Item(40, 9).Value = 1633.3334  # I40: was 1587.375
$ws.Cells.Item(40, 10).Value = 1987.375  # J40: was 2000
$ws.Cells.Item(40, 11).Value = 1633.3334  # K40: was 1587.375
$ws.Cells.Item(40, 12).Value = 1987.375  # L40: was 2000
$ws.Cells.Item(40, 13).Value = -1458.3334  # M40: was -1412.375
$ws.Cells.Item(40, 14).Value = -2337.375  # N40: was -2350
$ws.Cells.Item(55, 8).Value = 400  # H55: was 562.5
$ws.Cells.Item(55, 9).Value = 200.25  # I55: was 200.2
$ws.Cells.Item(55, 10).Value = 799.5  # J55: was 1166.3334
$ws.Cells.Item(55, 11).Value = 200.25  # K55: was 200.2
$ws.Cells.Item(55, 12).Value = 799.5  # L55: was 1166.3334
$ws.Cells.Item(55, 13).Value = 13.75  # M55: was 13.80000000000001
$ws.Cells.Item(55, 14).Value = -1227.5  # N55: was -1594.3334
$ws.Cells.Item(137, 8).Value = 1200.2667  # H137: was 1234.1395
$ws.Cells.Item(137, 9).Value = 1194.3715  # I137: was 1238.1515
$ws.Cells.Item(137, 11).Value = 3583.1145  # K137: was 3714.4545
$ws.Cells.Item(137, 13).Value = -1033.1145  # M137: was -1164.4545

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(35, 8).Value = 16509.25  # H35: was 13918.25
$ws.Cells.Item(35, 9).Value = 18012.334  # I35: was 13918.25
$ws.Cells.Item(35, 10).Value = 12000  # J35: was 0
$ws.Cells.Item(35, 11).Value = 18012.334  # K35: was 13918.25
$ws.Cells.Item(35, 12).Value = 12000  # L35: was 0
$ws.Cells.Item(35, 13).Value = -17606.334  # M35: was -13512.25
$ws.Cells.Item(35, 14).Value = -12812  # N35: was NEW
$ws.Cells.Item(45, 8).Value = 1003.6875  # H45: was 1173.6
$ws.Cells.Item(45, 9).Value = 942.2308  # I45: was 1132.2858
$ws.Cells.Item(45, 11).Value = 942.2308  # K45: was 1132.2858
$ws.Cells.Item(45, 13).Value = -565.2308  # M45: was -755.2858000000001
$ws.Cells.Item(122, 8).Value = 2615.4324  # H122: was 2855.4243
$ws.Cells.Item(122, 9).Value = 2545.5312  # I122: was 2889.4443
$ws.Cells.Item(122, 10).Value = 3062.8  # J122: was 2702.3333
$ws.Cells.Item(122, 11).Value = 7636.5936  # K122: was 8668.332900000001
$ws.Cells.Item(122, 12).Value = 9188.400000000001  # L122: was 8106.999899999999
$ws.Cells.Item(122, 13).Value = -5186.5936  # M122: was -6218.332900000001
$ws.Cells.Item(122, 14).Value = -14088.4  # N122: was -13006.9999
$ws.Cells.Item(131, 8).Value = 58358  # H131: was 59715
$ws.Cells.Item(131, 10).Value = 58358  # J131: was 59715
$ws.Cells.Item(131, 12).Value = 58358  # L131: was 59715
$ws.Cells.Item(131, 14).Value = -68438  # N131: was -69795
$ws.Cells.Item(132, 8).Value = 3438.7966  # H132: was 3901.5881
$ws.Cells.Item(132, 9).Value = 4099.7896  # I132: was 4674.364
$ws.Cells.Item(132, 10).Value = 2242.7144  # J132: was 2484.8333
$ws.Cells.Item(132, 11).Value = 12299.3688  # K132: was 14023.092
$ws.Cells.Item(132, 12).Value = 6728.1432  # L132: was 7454.499899999999
$ws.Cells.Item(132, 13).Value = -9769.3688  # M132: was -11493.092
$ws.Cells.Item(132, 14).Value = -11788.1432  # N132: was -12514.4999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(97, 8).Value = 38666.332  # H97: was 38499.668
$ws.Cells.Item(97, 10).Value = 38666.332  # J97: was 38499.668
$ws.Cells.Item(97, 12).Value = 38666.332  # L97: was 38499.668
$ws.Cells.Item(97, 14).Value = -40648.332  # N97: was -40481.668
$ws.Cells.Item(100, 8).Value = 2536630  # H100: was 10000000
$ws.Cells.Item(100, 10).Value = 2536630  # J100: was 10000000
$ws.Cells.Item(100, 12).Value = 2536630  # L100: was 10000000
$ws.Cells.Item(100, 14).Value = -2538794  # N100: was -10002164

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(15, 8).Value = 15  # H15: was 400
$ws.Cells.Item(15, 9).Value = 15  # I15: was 0
$ws.Cells.Item(15, 10).Value = 0  # J15: was 400
$ws.Cells.Item(15, 11).Value = 45  # K15: was 0
$ws.Cells.Item(15, 12).Value = 0  # L15: was 1200
$ws.Cells.Item(15, 14).ClearContents()  # N15: was -1480
$ws.Cells.Item(15, 13).Value = 95  # M15: was NEW
$ws.Cells.Item(94, 8).Value = 3126.6667  # H94: was 3380.9092
$ws.Cells.Item(94, 9).Value = 1950  # I94: was 2000
$ws.Cells.Item(94, 10).Value = 3462.8572  # J94: was 3519
$ws.Cells.Item(94, 11).Value = 5850  # K94: was 6000
$ws.Cells.Item(94, 12).Value = 10388.5716  # L94: was 10557
$ws.Cells.Item(94, 13).Value = -5174  # M94: was -5324
$ws.Cells.Item(94, 14).Value = -11740.5716  # N94: was -11909
$ws.Cells.Item(116, 8).Value = 1960.375  # H116: was 3569.8572
$ws.Cells.Item(116, 9).Value = 781.8333  # I116: was 1264.25
$ws.Cells.Item(116, 10).Value = 5496  # J116: was 6644
$ws.Cells.Item(116, 11).Value = 2345.4999  # K116: was 3792.75
$ws.Cells.Item(116, 12).Value = 16488  # L116: was 19932
$ws.Cells.Item(116, 13).Value = 1096.5001  # M116: was -350.75
$ws.Cells.Item(116, 14).Value = -23372  # N116: was -26816
$ws.Cells.Item(131, 8).Value = 10427090  # H131: was 10010297
$ws.Cells.Item(131, 10).Value = 10880423  # J131: was 10427375
$ws.Cells.Item(131, 12).Value = 32641269  # L131: was 31282125
$ws.Cells.Item(131, 14).Value = -32651349  # N131: was -31292205

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2299.484  # H122: was 2449.5334
$ws.Cells.Item(122, 9).Value = 1935.5652  # I122: was 1991.8182
$ws.Cells.Item(122, 10).Value = 3345.75  # J122: was 3708.25
$ws.Cells.Item(122, 11).Value = 5806.6956  # K122: was 5975.4546
$ws.Cells.Item(122, 12).Value = 10037.25  # L122: was 11124.75
$ws.Cells.Item(122, 13).Value = -3356.6956  # M122: was -3525.4546
$ws.Cells.Item(122, 14).Value = -14937.25  # N122: was -16024.75
$ws.Cells.Item(123, 8).Value = 8922.684999999999  # H123: was 8922.916999999999
$ws.Cells.Item(123, 10).Value = 8922.684999999999  # J123: was 8922.916999999999
$ws.Cells.Item(123, 12).Value = 8922.684999999999  # L123: was 8922.916999999999
$ws.Cells.Item(123, 14).Value = -13822.685  # N123: was -13822.917
$ws.Cells.Item(131, 8).Value = 38995.668  # H131: was 48659.332
$ws.Cells.Item(131, 10).Value = 38995.668  # J131: was 48659.332
$ws.Cells.Item(131, 12).Value = 38995.668  # L131: was 48659.332
$ws.Cells.Item(131, 14).Value = -49075.668  # N131: was -58739.332
$ws.Cells.Item(132, 8).Value = 2344.1614  # H132: was 2043.175
$ws.Cells.Item(132, 9).Value = 1493.6316  # I132: was 1374.9131
$ws.Cells.Item(132, 10).Value = 3690.8333  # J132: was 2947.2942
$ws.Cells.Item(132, 11).Value = 4480.8948  # K132: was 4124.7393
$ws.Cells.Item(132, 12).Value = 11072.4999  # L132: was 8841.882599999999
$ws.Cells.Item(132, 13).Value = -1950.8948  # M132: was -1594.7393
$ws.Cells.Item(132, 14).Value = -16132.4999  # N132: was -13901.8826

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value = 15026245  # H20: was 8592571
$ws.Cells.Item(20, 9).Value = 20033334  # I20: was 15028000
$ws.Cells.Item(20, 10).Value = 4980  # J20: was 12000
$ws.Cells.Item(20, 11).Value = 20033334  # K20: was 15028000
$ws.Cells.Item(20, 12).Value = 4980  # L20: was 12000
$ws.Cells.Item(20, 13).Value = -20033108  # M20: was -15027774
$ws.Cells.Item(20, 14).Value = -5432  # N20: was -12452
$ws.Cells.Item(22, 8).Value = 707.4783  # H22: was 763.1
$ws.Cells.Item(22, 9).Value = 716.5454999999999  # I22: was 808
$ws.Cells.Item(22, 10).Value = 699.1667  # J22: was 726.36365
$ws.Cells.Item(22, 11).Value = 716.5454999999999  # K22: was 808
$ws.Cells.Item(22, 12).Value = 699.1667  # L22: was 726.36365
$ws.Cells.Item(22, 13).Value = -421.5454999999999  # M22: was -513
$ws.Cells.Item(22, 14).Value = -1289.1667  # N22: was -1316.36365
$ws.Cells.Item(27, 8).Value = 707.4783  # H27: was 763.1
$ws.Cells.Item(27, 9).Value = 716.5454999999999  # I27: was 808
$ws.Cells.Item(27, 10).Value = 699.1667  # J27: was 726.36365
$ws.Cells.Item(27, 11).Value = 716.5454999999999  # K27: was 808
$ws.Cells.Item(27, 12).Value = 699.1667  # L27: was 726.36365
$ws.Cells.Item(27, 13).Value = -609.5454999999999  # M27: was -701
$ws.Cells.Item(27, 14).Value = -913.1667  # N27: was -940.36365
$ws.Cells.Item(46, 8).Value = 837.7838  # H46: was 1025.0625
$ws.Cells.Item(46, 9).Value = 621.2121  # I46: was 781.9091
$ws.Cells.Item(46, 10).Value = 2624.5  # J46: was 1560
$ws.Cells.Item(46, 11).Value = 621.2121  # K46: was 781.9091
$ws.Cells.Item(46, 12).Value = 2624.5  # L46: was 1560
$ws.Cells.Item(46, 13).Value = -433.2121  # M46: was -593.9091
$ws.Cells.Item(46, 14).Value = -3000.5  # N46: was -1936
$ws.Cells.Item(61, 8).Value = 3500  # H61: was 3200
$ws.Cells.Item(61, 9).Value = 2500  # I61: was 1650
$ws.Cells.Item(61, 10).Value = 4500  # J61: was 4750
$ws.Cells.Item(61, 11).Value = 2500  # K61: was 1650
$ws.Cells.Item(61, 12).Value = 4500  # L61: was 4750
$ws.Cells.Item(61, 13).Value = -2298  # M61: was -1448
$ws.Cells.Item(61, 14).Value = -4904  # N61: was -5154
$ws.Cells.Item(68, 8).Value = 3536.6365  # H68: was 3167.8235
$ws.Cells.Item(68, 9).Value = 2862.5  # I68: was 2631.818
$ws.Cells.Item(68, 10).Value = 5334.3335  # J68: was 4150.5
$ws.Cells.Item(68, 11).Value = 2862.5  # K68: was 2631.818
$ws.Cells.Item(68, 12).Value = 5334.3335  # L68: was 4150.5
$ws.Cells.Item(68, 13).Value = -2113.5  # M68: was -1882.818
$ws.Cells.Item(68, 14).Value = -6832.3335  # N68: was -5648.5
$ws.Cells.Item(71, 8).Value = 3536.6365  # H71: was 3167.8235
$ws.Cells.Item(71, 9).Value = 2862.5  # I71: was 2631.818
$ws.Cells.Item(71, 10).Value = 5334.3335  # J71: was 4150.5
$ws.Cells.Item(71, 11).Value = 14312.5  # K71: was 13159.09
$ws.Cells.Item(71, 12).Value = 26671.6675  # L71: was 20752.5
$ws.Cells.Item(71, 13).Value = -10568.5  # M71: was -9415.09
$ws.Cells.Item(71, 14).Value = -34159.6675  # N71: was -28240.5
$ws.Cells.Item(113, 8).Value = 3500  # H113: was 3200
$ws.Cells.Item(113, 9).Value = 2500  # I113: was 1650
$ws.Cells.Item(113, 10).Value = 4500  # J113: was 4750
$ws.Cells.Item(113, 11).Value = 2500  # K113: was 1650
$ws.Cells.Item(113, 12).Value = 4500  # L113: was 4750
$ws.Cells.Item(113, 13).Value = -330  # M113: was 520
$ws.Cells.Item(113, 14).Value = -8840  # N113: was -9090
$ws.Cells.Item(131, 8).Value = 30326  # H131: was 0
$ws.Cells.Item(131, 10).Value = 30326  # J131: was 0
$ws.Cells.Item(131, 12).Value = 30326  # L131: was 0
$ws.Cells.Item(131, 14).Value = -40406  # N131: was NEW
$ws.Cells.Item(132, 8).Value = 4671.2856  # H132: was 3241
$ws.Cells.Item(132, 9).Value = 4737.5  # I132: was 2846.0952
$ws.Cells.Item(132, 10).Value = 4583  # J132: was 4899.6
$ws.Cells.Item(132, 11).Value = 14212.5  # K132: was 8538.285600000001
$ws.Cells.Item(132, 12).Value = 13749  # L132: was 14698.8
$ws.Cells.Item(132, 13).Value = -11682.5  # M132: was -6008.285600000001
$ws.Cells.Item(132, 14).Value = -18809  # N132: was -19758.8
$ws.Cells.Item(133, 8).Value = 55662  # H133: was 84660
$ws.Cells.Item(133, 10).Value = 55662  # J133: was 84660
$ws.Cells.Item(133, 12).Value = 55662  # L133: was 84660
$ws.Cells.Item(133, 14).Value = -60722  # N133: was -89720

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 8).Value = 19503.166  # H123: was 20377.678
$ws.Cells.Item(123, 10).Value = 19503.166  # J123: was 20377.678
$ws.Cells.Item(123, 12).Value = 19503.166  # L123: was 20377.678
$ws.Cells.Item(123, 14).Value = -29303.166  # N123: was -30177.678
$ws.Cells.Item(132, 8).Value = 1946.9736  # H132: was 1781.0465
$ws.Cells.Item(132, 9).Value = 1299.6364  # I132: was 1203.68
$ws.Cells.Item(132, 10).Value = 2837.0625  # J132: was 2582.9443
$ws.Cells.Item(132, 11).Value = 3898.9092  # K132: was 3611.04
$ws.Cells.Item(132, 12).Value = 8511.1875  # L132: was 7748.8329
$ws.Cells.Item(132, 13).Value = -1368.9092  # M132: was -1081.04
$ws.Cells.Item(132, 14).Value = -13571.1875  # N132: was -12808.8329
